$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3384
$ws1.Range("F5").Value = 6959
$ws1.Range("F6").Value = 2396
$ws1.Range("F8").Value = 106
$ws1.Range("F10").Value = 37

# Sheet "全部类型" (sheet4): update the same events' "想去人数" (column F) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3384
$ws4.Range("F6").Value = 6959
$ws4.Range("F7").Value = 2396
$ws4.Range("F9").Value = 106
$ws4.Range("F11").Value = 37

$wb.Save()
